# "cambios de las vistas" - view/filter changes on Hoja1 (+ Hoja2 tab deselect)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# Make Hoja1 the active/selected tab (was Hoja2) and zoom it to 160%.
$ws1.Activate()
$excel.ActiveWindow.Zoom = 160

# Move the selection on Hoja1 to A5.
$ws1.Range("A5").Select()

# Apply an AutoFilter on column A (LOCAL) of the data range A1:E23,
# keeping only MP100 / MP108 / MP15 / MP53 visible (rest of the rows hide).
$rng1 = $ws1.Range("A1:E23")
$rng1.AutoFilter()
$rng1.AutoFilter(1, @("MP100", "MP108", "MP15", "MP53"), 7)

# Keep the _FilterDatabase defined name in sync with the new filter range.
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -eq "Hoja1!_FilterDatabase") {
        $n.RefersTo = "=Hoja1!`$A`$1:`$E`$23"
    }
}

# Hoja2 is no longer the selected tab; restore its own selection.
$ws2.Activate()
$ws2.Range("F11").Select()

$ws1.Activate()
